$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 172 - this shifts the existing rows 172:242 down to 173:243,
# carrying their values/formatting with them (matches the target diff, where every
# row from 173 to 243 now holds what used to be one row above it).
$ws.Rows(172).Insert()

# Populate the newly inserted row 172 with a new data record. Columns that are
# identical to the (now shifted-down) neighbouring record are simply copied from
# row 173; the remaining columns get the new values from the diff.
$ws.Range("A172").Value = 10
$ws.Range("B172").Value = "Vega Modelo de Temuco"
$ws.Range("C172").Value = "La Araucanía"
$ws.Range("D172").Value = 44510
$ws.Range("E172").Value = 9
$ws.Range("F172").Value = "Fruta"
$ws.Range("G172").Value = 100108
$ws.Range("H172").Value = "Tropicales y subtropicales"
$ws.Range("I172").Value = 100108002
$ws.Range("J172").Value = "Mango"
$ws.Range("K172").Value = "Sin especificar"
$ws.Range("L172").Value = "Primera"
$ws.Range("M172").Value = 55
$ws.Range("N172").Value = 8000
$ws.Range("O172").Value = 8000
$ws.Range("P172").Value = 8000
$ws.Range("Q172").Value = "$/bandeja 4 kilos"
$ws.Range("R172").Value = "Perú"
$ws.Range("S172").Value = 2000
$ws.Range("T172").Value = 4
